# Apply updated crypto price/volume figures and the ZBToken/CoinExToken row swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'279.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.47%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.95%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.833"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.83%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'0.20%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.044"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.98%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.304"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.82%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8951"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.62%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1543"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.10%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06658"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'29.64%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07490"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.68%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02948"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.43%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08999"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.12%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001569"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.36%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.68%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.21%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'0.66%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.328"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.27%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.228"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.93%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3145"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.85%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'1.07%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.900"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.19%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04408"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.29%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").Value = "'0.1504"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'8.97%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.35%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004276"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'10.08%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E28").Value = "'-1.60%"
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'-14.57%"
$ws.Range("E29").Style = "Normal"
$ws.Range("E40").Value = "'-2.08%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006617"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-3.38%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1402"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'18.94%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002060"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.05%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01100"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.95%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005550"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.17%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E47").Value = "'-8.55%"
$ws.Range("E47").Style = "Normal"
